$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 4; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 8; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 10; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 16; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 20; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 25; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 27; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 34; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 37; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 49; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 52; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 57; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 60; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 68; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 71; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 75; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 83; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 87; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 91; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 94; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 114; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 115; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 121; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 124; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 125; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 137; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 139; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 141; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 143; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 145; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 155; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 177; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 189; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 200; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 201; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 205; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 216; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 217; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 224; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 234; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 235; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 236; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 239; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 246; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 249; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 259; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 268; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 288; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 290; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 302; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 309; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 310; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 311; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 314; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 324; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 326; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 333; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 368; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 369; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 374; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 388; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 407; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 417; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 422; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 434; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 445; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 447; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 448; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 462; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 469; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 478; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 481; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 482; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 483; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 484; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 485; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 486; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 493; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 495; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 497; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 504; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 506; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 508; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 511; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 518; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 522; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 525; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 529; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}

Write-Output "Updated $($updates.Count) rows"